$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column D ("Terms Typically Offered"),
# shifting that column to G. This reflects the new requirement separation:
# Corequisites, Concurrent, Recommended.
$ws.Range("D1:F1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Fill "NA" for the new columns across all data rows (2-27)
$ws.Range("D2:F27").Value = "NA"
